$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.216.96"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").Value = "'2.061.55"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'230.21"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("D7").Value = "'60.55"
$ws.Range("E7").Value = "  +9.68%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("D10").Value = "'0.0811"
$ws.Range("E10").Value = "  +4.47%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "'14.81"
$ws.Range("E12").Value = "  +5.45%  "
$ws.Range("D13").Value = "'2.364.12"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "'21.38"
$ws.Range("E14").Value = "  +8.25%  "
$ws.Range("D15").Value = "'0.761"
$ws.Range("E15").Value = "  +3.62%  "
$ws.Range("D16").Value = "'5.31"
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "'2.060.08"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").Value = "'38.084.29"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("E19").Value = "  +2.24%  "
$ws.Range("D20").Value = "'70.06"
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").Value = "'0.0₃0836"
$ws.Range("E21").Value = "  +3.11%  "
$ws.Range("D22").Value = "'225.86"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("D23").Value = "'0.995"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'2.25"
$ws.Range("E25").Value = "  +3.49%  "
$ws.Range("D26").Value = "'9.32"
$ws.Range("E26").Value = "  +4.30%  "
$ws.Range("D27").Value = "'166.49"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E28").Value = "  +6.27%  "
$ws.Range("D29").Value = "'19.11"
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("D32").Value = "'4.56"
$ws.Range("E32").Value = "  +3.69%  "
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("E34").Value = "  +9.80%  "
$ws.Range("D35").Value = "'0.0607"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.31"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'6.21"
$ws.Range("E37").Value = "  +16.18%  "
$ws.Range("D38").Value = "'3.31"
$ws.Range("E38").Value = "  +5.41%  "
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'1.526.93"
$ws.Range("E40").Value = "  +4.33%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'98.19"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'17.12"
$ws.Range("E42").Value = "  +7.73%  "
$ws.Range("D43").Value = "'0.0217"
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("D45").Value = "'0.0925"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("D47").Value = "'4.07"
$ws.Range("E47").Value = "  -5.21%  "
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("D50").Value = "'7.14"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "'2.251.65"
$ws.Range("E51").Value = "  +2.45%  "
